# Edit script: insert 4 new weekly rows of Mandarina price data
# (Fruta / hortaliza, semanal) right before the existing row 361,
# pushing the former rows 361-378 down to 365-382.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 361..364 (existing data shifts down to 365..382)
$ws.Rows("361:364").Insert()

# Common column values shared by all rows in this data block
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"

# New row 361: Clemenuless / Primera
$r = 361
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 45147
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = "Clemenuless"
$ws.Cells.Item($r, 12).Value2 = "Primera"
$ws.Cells.Item($r, 13).Value2 = 100
$ws.Cells.Item($r, 14).Value2 = 8000
$ws.Cells.Item($r, 15).Value2 = 8000
$ws.Cells.Item($r, 16).Value2 = 8000
$ws.Cells.Item($r, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item($r, 19).Value2 = 800
$ws.Cells.Item($r, 20).Value2 = 10

# New row 362: Clemenuless / Segunda
$r = 362
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 45147
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = "Clemenuless"
$ws.Cells.Item($r, 12).Value2 = "Segunda"
$ws.Cells.Item($r, 13).Value2 = 60
$ws.Cells.Item($r, 14).Value2 = 6000
$ws.Cells.Item($r, 15).Value2 = 6000
$ws.Cells.Item($r, 16).Value2 = 6000
$ws.Cells.Item($r, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item($r, 19).Value2 = 600
$ws.Cells.Item($r, 20).Value2 = 10

# New row 363: Murcott / Primera
$r = 363
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 45147
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = "Murcott"
$ws.Cells.Item($r, 12).Value2 = "Primera"
$ws.Cells.Item($r, 13).Value2 = 60
$ws.Cells.Item($r, 14).Value2 = 11000
$ws.Cells.Item($r, 15).Value2 = 11000
$ws.Cells.Item($r, 16).Value2 = 11000
$ws.Cells.Item($r, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item($r, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value2 = 611
$ws.Cells.Item($r, 20).Value2 = 18

# New row 364: Murcott / Segunda
$r = 364
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 45147
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = "Murcott"
$ws.Cells.Item($r, 12).Value2 = "Segunda"
$ws.Cells.Item($r, 13).Value2 = 60
$ws.Cells.Item($r, 14).Value2 = 9000
$ws.Cells.Item($r, 15).Value2 = 9000
$ws.Cells.Item($r, 16).Value2 = 9000
$ws.Cells.Item($r, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item($r, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value2 = 500
$ws.Cells.Item($r, 20).Value2 = 18

# Ensure the date cells (column D) for the new rows use the same date/time
# number format as the rest of the column (style already carried over from
# the insert, but set explicitly to be safe).
$ws.Range("D361:D364").NumberFormat = $ws.Range("D365").NumberFormat
